$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns: S1 = "Idade ignorada", T1 = "Total"
$ws.Range("S1").Value = "Idade ignorada"
$ws.Range("T1").Value = "Total"

# Existing rows 2-6 get a new "Total" value in column T (S left blank, but touched
# so the cell exists in the sheet - matches the target XML's empty <c r="S2"/> etc.)
$ws.Range("S2").Style = "Normal"
$ws.Range("T2").Value = 2050
$ws.Range("S3").Style = "Normal"
$ws.Range("T3").Value = 223
$ws.Range("S4").Style = "Normal"
$ws.Range("T4").Value = 1082
$ws.Range("S5").Style = "Normal"
$ws.Range("T5").Value = 236
$ws.Range("S6").Style = "Normal"
$ws.Range("T6").Value = 1358

# New row 7: "Outros" (Others)
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 147
$ws.Range("C7").Value = 7
$ws.Range("D7").Value = 19
$ws.Range("E7").Value = 42
$ws.Range("F7").Value = 99
$ws.Range("G7").Value = 87
$ws.Range("H7").Value = 84
$ws.Range("I7").Value = 107
$ws.Range("J7").Value = 105
$ws.Range("K7").Value = 104
$ws.Range("L7").Value = 105
$ws.Range("M7").Value = 108
$ws.Range("N7").Value = 99
$ws.Range("O7").Value = 106
$ws.Range("P7").Value = 124
$ws.Range("Q7").Value = 132
$ws.Range("R7").Value = 384
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 1860

# New row 8: "Total" (grand total)
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 169
$ws.Range("C8").Value = 13
$ws.Range("D8").Value = 24
$ws.Range("E8").Value = 51
$ws.Range("F8").Value = 115
$ws.Range("G8").Value = 109
$ws.Range("H8").Value = 115
$ws.Range("I8").Value = 169
$ws.Range("J8").Value = 188
$ws.Range("K8").Value = 285
$ws.Range("L8").Value = 367
$ws.Range("M8").Value = 468
$ws.Range("N8").Value = 515
$ws.Range("O8").Value = 616
$ws.Range("P8").Value = 623
$ws.Range("Q8").Value = 828
$ws.Range("R8").Value = 2153
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 6809
